# ArcaDigital WebApp - Upload/Extract test data edit.
# Renames the sample product (row 2) so it reads as an auto-upload test
# record, and aligns the barcode column with the internal product code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nombre (A2), Descripcion (P2) and Nombre secundario (Q2) all previously
# held the literal product name "BILLETERA D&G COLOR NEGRO" - replace it
# with the new auto-upload test marker.
$ws.Range("A2").Value = "Auto Upload Test File"
$ws.Range("P2").Value = "Auto Upload Test File"
$ws.Range("Q2").Value = "Auto Upload Test File"

# Cod barras (T2) now mirrors the internal code (Codigo Interno, B2).
$ws.Range("T2").Value = "BI001"

# Leave the cursor on E13, matching the state the file was saved in.
$ws.Range("E13").Select()
